$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.340.31'
$ws.Range("E2").Value = '  -3.09%  '

$ws.Range("D3").Value = '1.831.54'
$ws.Range("E3").Value = '  -2.61%  '

$orig = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = $orig
$ws.Range("E4").Value = '  -0.05%  '

$orig = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.26'
$ws.Range("D5").Style = $orig
$ws.Range("E5").Value = '  -7.82%  '

$ws.Range("E6").Value = '  -0.03%  '

$orig = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5194'
$ws.Range("D7").Style = $orig
$ws.Range("E7").Value = '  -1.89%  '

$orig = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3215'
$ws.Range("D8").Style = $orig
$ws.Range("E8").Value = '  -9.06%  '

$orig = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06729'
$ws.Range("D9").Style = $orig
$ws.Range("E9").Value = '  -4.41%  '

$orig = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.66'
$ws.Range("D10").Style = $orig

$orig = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7648'
$ws.Range("D11").Style = $orig
$ws.Range("E11").Value = '  -6.94%  '

$ws.Range("D12").Value = '1.880.27'
$ws.Range("E12").Value = '  -0.30%  '

$orig = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07673'
$ws.Range("D13").Style = $orig
$ws.Range("E13").Value = '  -1.85%  '

$orig = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.74'
$ws.Range("D14").Style = $orig
$ws.Range("E14").Value = '  -2.27%  '

$orig = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.017'
$ws.Range("D15").Style = $orig
$ws.Range("E15").Value = '  -3.61%  '

$ws.Range("E16").Value = '  -0.01%  '

$orig = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.03'
$ws.Range("D17").Style = $orig
$ws.Range("E17").Value = '  -4.05%  '

$orig = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9999'
$ws.Range("D18").Style = $orig
$ws.Range("E18").Value = '  +0.03%  '

$orig = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007895'
$ws.Range("D19").Style = $orig
$ws.Range("E19").Value = '  -3.34%  '

$ws.Range("D20").Value = '26.392.66'
$ws.Range("E20").Value = '  -3.10%  '

$ws.Range("D21").Value = '2.088.07'
$ws.Range("E21").Value = '  -2.09%  '

$orig = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.528'
$ws.Range("D22").Style = $orig
$ws.Range("E22").Value = '  -5.15%  '

$orig = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.413'
$ws.Range("D23").Style = $orig
$ws.Range("E23").Value = '  -7.34%  '

$orig = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.907'
$ws.Range("D24").Style = $orig
$ws.Range("E24").Value = '  -5.54%  '

$orig = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.279'
$ws.Range("D25").Style = $orig
$ws.Range("E25").Value = '  -5.01%  '

$orig = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.06'
$ws.Range("D26").Style = $orig
$ws.Range("E26").Value = '  -1.42%  '

$orig = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.639'
$ws.Range("D27").Style = $orig
$ws.Range("E27").Value = '  -2.14%  '

$orig = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.89'
$ws.Range("D28").Style = $orig
$ws.Range("E28").Value = '  -4.11%  '

$orig = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.10'
$ws.Range("D29").Style = $orig
$ws.Range("E29").Value = '  -3.20%  '

$ws.Range("E30").Value = '  -5.17%  '

$orig = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.129'
$ws.Range("D31").Style = $orig
$ws.Range("E31").Value = '  -5.90%  '

$orig = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08738'
$ws.Range("D32").Style = $orig
$ws.Range("E32").Value = '  -2.41%  '

$orig = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04844'
$ws.Range("D33").Style = $orig
$ws.Range("E33").Value = '  -2.02%  '

$ws.Range("E34").Value = '  -4.87%  '

$orig = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.843'
$ws.Range("D35").Style = $orig
$ws.Range("E35").Value = '  -2.30%  '

$orig = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6821'
$ws.Range("D36").Style = $orig
$ws.Range("E36").Value = '  -8.98%  '

$orig = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.090'
$ws.Range("D37").Style = $orig
$ws.Range("E37").Value = '  -6.74%  '

$orig = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01774'
$ws.Range("D38").Style = $orig
$ws.Range("E38").Value = '  -5.97%  '

$orig = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.222'
$ws.Range("D39").Style = $orig
$ws.Range("E39").Value = '  -8.37%  '

$orig = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4897'
$ws.Range("D40").Style = $orig
$ws.Range("E40").Value = '  -7.98%  '

$orig = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '112.51'
$ws.Range("D41").Style = $orig
$ws.Range("E41").Value = '  -3.83%  '

$orig = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8878'
$ws.Range("D42").Style = $orig
$ws.Range("E42").Value = '  -8.79%  '

$orig = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.116'
$ws.Range("D43").Style = $orig
$ws.Range("E43").Value = '  -3.31%  '

$orig = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("D44").Style = $orig
$ws.Range("E44").Value = '  -0.02%  '

$orig = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.691'
$ws.Range("D45").Style = $orig

$orig = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4195'
$ws.Range("D46").Style = $orig
$ws.Range("E46").Value = '  -9.06%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$orig = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.086'
$ws.Range("D47").Style = $orig
$ws.Range("E47").Value = '  -4.58%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$orig = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1253'
$ws.Range("D48").Style = $orig
$ws.Range("E48").Value = '  -8.74%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$orig = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05873'
$ws.Range("D49").Style = $orig
$ws.Range("E49").Value = '  -1.41%  '

$orig = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.29'
$ws.Range("D50").Style = $orig
$ws.Range("E50").Value = '  -3.90%  '

$orig = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.28'
$ws.Range("D51").Style = $orig
$ws.Range("E51").Value = '  -4.24%  '
